# Weekly update: insert the newest week's "Apio" (Vega Monumental Concepción)
# price record as a new row 182, pushing all the older weekly records down by
# one row (dimension grows from R302 to R303).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 182; Excel shifts rows 182..302 down to
# 183..303 and the used range / dimension grows to A1:R303 automatically.
$ws.Rows.Item(182).Insert()

# Populate the newly inserted row with this week's data.
$ws.Cells.Item(182, 1).Value = 11
$ws.Cells.Item(182, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(182, 3).Value = "Bíobío"
$ws.Cells.Item(182, 4).Value = 44762
$ws.Cells.Item(182, 5).Value = 8
$ws.Cells.Item(182, 6).Value = 100112017
$ws.Cells.Item(182, 7).Value = "Apio"
$ws.Cells.Item(182, 8).Value = "Americana (o)"
$ws.Cells.Item(182, 9).Value = "Primera"
$ws.Cells.Item(182, 10).Value = 200
$ws.Cells.Item(182, 11).Value = 8000
$ws.Cells.Item(182, 12).Value = 8500
$ws.Cells.Item(182, 13).Value = 8250
$ws.Cells.Item(182, 14).Value = "$/docena de matas"
$ws.Cells.Item(182, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(182, 16).Value = 1375
$ws.Cells.Item(182, 17).Value = 6
$ws.Cells.Item(182, 18).Value = "Hortaliza"
